$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 3
$ws.Range("F2").Value = 109
$ws.Range("H2").Value = "living_rooms"
$ws.Range("I2").Value = "distractor"
$ws.Range("K2").Value = "f"
$ws.Range("L2").Value = "stimuli/img_bbs77.png"
$ws.Range("M2").Value = 31.64444444444445
$ws.Range("N2").Value = 21.26666666666667
$ws.Range("O2").Value = 26.45555555555556
$ws.Range("P2").Value = 45
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 2

# Row 3
$ws.Range("C3").Value = 3
$ws.Range("F3").Value = 110
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "distractor"
$ws.Range("K3").Value = "f"
$ws.Range("L3").Value = "stimuli/img_pt3d7.png"
$ws.Range("M3").Value = 65.08571428571429
$ws.Range("N3").Value = 44.65714285714286
$ws.Range("O3").Value = 54.87142857142857
$ws.Range("P3").Value = 35
$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 4
$ws.Range("S3").Value = 4

# Row 4
$ws.Range("C4").Value = 3
$ws.Range("F4").Value = 111
$ws.Range("H4").Value = "living_rooms"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_xy930.png"
$ws.Range("M4").Value = 70.5952380952381
$ws.Range("N4").Value = 49.47619047619047
$ws.Range("O4").Value = 60.03571428571429
$ws.Range("P4").Value = 42
$ws.Range("Q4").Value = 6
$ws.Range("R4").Value = 6
$ws.Range("S4").Value = 6

# Row 5
$ws.Range("C5").Value = 3
$ws.Range("F5").Value = 112
$ws.Range("H5").Value = "living_rooms"
$ws.Range("I5").Value = "target"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_abobq.png"
$ws.Range("M5").Value = 75.1842105263158
$ws.Range("N5").Value = 54.13157894736842
$ws.Range("O5").Value = 64.65789473684211
$ws.Range("P5").Value = 38
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 6
$ws.Range("S5").Value = 6

# Row 6
$ws.Range("C6").Value = 3
$ws.Range("F6").Value = 113
$ws.Range("H6").Value = "living_rooms"
$ws.Range("I6").Value = "target"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_di6f0.png"
$ws.Range("M6").Value = 94.04347826086956
$ws.Range("N6").Value = 83.34782608695652
$ws.Range("O6").Value = 88.69565217391303
$ws.Range("P6").Value = 46
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 10
$ws.Range("S6").Value = 10

# Row 7
$ws.Range("C7").Value = 3
$ws.Range("F7").Value = 114
$ws.Range("H7").Value = "living_rooms"
$ws.Range("I7").Value = "target"
$ws.Range("K7").Value = "j"
$ws.Range("L7").Value = "stimuli/img_eh0no.png"
$ws.Range("M7").Value = 53.66666666666666
$ws.Range("N7").Value = 36.02564102564103
$ws.Range("O7").Value = 44.84615384615385
$ws.Range("P7").Value = 39
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 3
$ws.Range("S7").Value = 3

# Row 8
$ws.Range("C8").Value = 3
$ws.Range("F8").Value = 115
$ws.Range("H8").Value = "living_rooms"
$ws.Range("I8").Value = "target"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_tujn3.png"
$ws.Range("M8").Value = 81.4090909090909
$ws.Range("N8").Value = 62.52272727272727
$ws.Range("O8").Value = 71.9659090909091
$ws.Range("P8").Value = 44
$ws.Range("Q8").Value = 8
$ws.Range("R8").Value = 8
$ws.Range("S8").Value = 8

# Row 9
$ws.Range("C9").Value = 3
$ws.Range("F9").Value = 116
$ws.Range("H9").Value = "living_rooms"
$ws.Range("I9").Value = "target"
$ws.Range("K9").Value = "j"
$ws.Range("L9").Value = "stimuli/img_4o8l0.png"
$ws.Range("M9").Value = 46.02173913043478
$ws.Range("N9").Value = 31.45652173913043
$ws.Range("O9").Value = 38.73913043478261
$ws.Range("P9").Value = 46
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 3
$ws.Range("S9").Value = 3

# Row 10
$ws.Range("C10").Value = 3
$ws.Range("F10").Value = 117
$ws.Range("H10").Value = "living_rooms"
$ws.Range("I10").Value = "target"
$ws.Range("K10").Value = "j"
$ws.Range("L10").Value = "stimuli/img_w8yhd.png"
$ws.Range("M10").Value = 55.74418604651163
$ws.Range("N10").Value = 38.90697674418605
$ws.Range("O10").Value = 47.32558139534883
$ws.Range("P10").Value = 43
$ws.Range("Q10").Value = 4
$ws.Range("R10").Value = 4
$ws.Range("S10").Value = 4

# Row 11
$ws.Range("C11").Value = 3
$ws.Range("F11").Value = 118
$ws.Range("H11").Value = "living_rooms"
$ws.Range("I11").Value = "target"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_6a0hu.png"
$ws.Range("M11").Value = 61.275
$ws.Range("N11").Value = 42.025
$ws.Range("O11").Value = 51.65
$ws.Range("P11").Value = 40
$ws.Range("Q11").Value = 4
$ws.Range("R11").Value = 4
$ws.Range("S11").Value = 4

# Row 12
$ws.Range("C12").Value = 3
$ws.Range("F12").Value = 119
$ws.Range("H12").Value = "kitchens"
$ws.Range("I12").Value = "distractor"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_gztbt.png"
$ws.Range("M12").Value = 55.06451612903226
$ws.Range("N12").Value = 26.09677419354839
$ws.Range("O12").Value = 40.58064516129032
$ws.Range("P12").Value = 31
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 2

# Row 13
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 120
$ws.Range("H13").Value = "living_rooms"
$ws.Range("I13").Value = "target"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_pey7u.png"
$ws.Range("M13").Value = 30.34883720930232
$ws.Range("N13").Value = 20.34883720930232
$ws.Range("O13").Value = 25.34883720930232
$ws.Range("P13").Value = 43
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 2
$ws.Range("S13").Value = 2

# Row 14
$ws.Range("C14").Value = 3
$ws.Range("F14").Value = 121
$ws.Range("H14").Value = "kitchens"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_7wquy.png"
$ws.Range("M14").Value = 50.59375
$ws.Range("N14").Value = 30.59375
$ws.Range("O14").Value = 40.59375
$ws.Range("P14").Value = 32
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 2
$ws.Range("S14").Value = 2

# Row 15
$ws.Range("C15").Value = 3
$ws.Range("F15").Value = 122
$ws.Range("H15").Value = "kitchens"
$ws.Range("I15").Value = "target"
$ws.Range("K15").Value = "j"
$ws.Range("L15").Value = "stimuli/img_xu1p3.png"
$ws.Range("M15").Value = 75.27659574468085
$ws.Range("N15").Value = 56.68085106382978
$ws.Range("O15").Value = 65.97872340425532
$ws.Range("P15").Value = 47
$ws.Range("Q15").Value = 7
$ws.Range("R15").Value = 7
$ws.Range("S15").Value = 7

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 123
$ws.Range("H16").Value = "bedrooms"
$ws.Range("I16").Value = "distractor"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_mdpr4.png"
$ws.Range("M16").Value = 74.04255319148936
$ws.Range("N16").Value = 54.70212765957447
$ws.Range("O16").Value = 64.37234042553192
$ws.Range("P16").Value = 47
$ws.Range("Q16").Value = 6
$ws.Range("R16").Value = 6
$ws.Range("S16").Value = 6

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("F17").Value = 124
$ws.Range("H17").Value = "kitchens"
$ws.Range("I17").Value = "target"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_die1d.png"
$ws.Range("M17").Value = 75.42857142857143
$ws.Range("N17").Value = 53.30952380952381
$ws.Range("O17").Value = 64.36904761904762
$ws.Range("P17").Value = 42
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = 6
$ws.Range("S17").Value = 6

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 125
$ws.Range("H18").Value = "living_rooms"
$ws.Range("I18").Value = "target"
$ws.Range("K18").Value = "j"
$ws.Range("L18").Value = "stimuli/img_cehin.png"
$ws.Range("M18").Value = 78.86363636363636
$ws.Range("N18").Value = 60.02272727272727
$ws.Range("O18").Value = 69.44318181818181
$ws.Range("P18").Value = 44
$ws.Range("Q18").Value = 7
$ws.Range("R18").Value = 7
$ws.Range("S18").Value = 7

# Row 19
$ws.Range("C19").Value = 3
$ws.Range("F19").Value = 126
$ws.Range("H19").Value = "living_rooms"
$ws.Range("I19").Value = "target"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_wz6x5.png"
$ws.Range("M19").Value = 68.3695652173913
$ws.Range("N19").Value = 48.47826086956522
$ws.Range("O19").Value = 58.42391304347826
$ws.Range("P19").Value = 46
$ws.Range("Q19").Value = 5
$ws.Range("R19").Value = 5
$ws.Range("S19").Value = 5

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("F20").Value = 127
$ws.Range("H20").Value = "living_rooms"
$ws.Range("I20").Value = "target"
$ws.Range("K20").Value = "j"
$ws.Range("L20").Value = "stimuli/img_6zz63.png"
$ws.Range("M20").Value = 87.66666666666667
$ws.Range("N20").Value = 70.6
$ws.Range("O20").Value = 79.13333333333333
$ws.Range("P20").Value = 45
$ws.Range("Q20").Value = 9
$ws.Range("R20").Value = 10
$ws.Range("S20").Value = 10

# Row 21
$ws.Range("C21").Value = 3
$ws.Range("F21").Value = 128
$ws.Range("H21").Value = "living_rooms"
$ws.Range("I21").Value = "target"
$ws.Range("K21").Value = "j"
$ws.Range("L21").Value = "stimuli/img_bj99b.png"
$ws.Range("M21").Value = 82.79069767441861
$ws.Range("N21").Value = 65.46511627906976
$ws.Range("O21").Value = 74.12790697674419
$ws.Range("P21").Value = 43
$ws.Range("Q21").Value = 8
$ws.Range("R21").Value = 8
$ws.Range("S21").Value = 8

# Row 22
$ws.Range("C22").Value = 3
$ws.Range("F22").Value = 129
$ws.Range("H22").Value = "living_rooms"
$ws.Range("I22").Value = "target"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_wgkqa.png"
$ws.Range("M22").Value = 87.25581395348837
$ws.Range("N22").Value = 71.13953488372093
$ws.Range("O22").Value = 79.19767441860465
$ws.Range("P22").Value = 43
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = 10
$ws.Range("S22").Value = 10

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 130
$ws.Range("H23").Value = "bedrooms"
$ws.Range("I23").Value = "distractor"
$ws.Range("K23").Value = "f"
$ws.Range("L23").Value = "stimuli/img_5p2ql.png"
$ws.Range("M23").Value = 89.19565217391305
$ws.Range("N23").Value = 72.52173913043478
$ws.Range("O23").Value = 80.8586956521739
$ws.Range("P23").Value = 46
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = 10
$ws.Range("S23").Value = 10

# Row 24
$ws.Range("C24").Value = 3
$ws.Range("F24").Value = 131
$ws.Range("H24").Value = "kitchens"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_kost0.png"
$ws.Range("M24").Value = 63.09090909090909
$ws.Range("N24").Value = 42.77272727272727
$ws.Range("O24").Value = 52.93181818181819
$ws.Range("P24").Value = 44
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("F25").Value = 132
$ws.Range("H25").Value = "living_rooms"
$ws.Range("I25").Value = "target"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_0kqc0.png"
$ws.Range("M25").Value = 43.74468085106383
$ws.Range("N25").Value = 27.14893617021277
$ws.Range("O25").Value = 35.4468085106383
$ws.Range("P25").Value = 47
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 2

# Row 26
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 133
$ws.Range("H26").Value = "living_rooms"
$ws.Range("I26").Value = "target"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_xbtev.png"
$ws.Range("M26").Value = 13.68181818181818
$ws.Range("N26").Value = 8.568181818181818
$ws.Range("O26").Value = 11.125
$ws.Range("P26").Value = 44
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 1

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("F27").Value = 134
$ws.Range("H27").Value = "living_rooms"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_95hiv.png"
$ws.Range("M27").Value = 84.04545454545455
$ws.Range("N27").Value = 67.31818181818181
$ws.Range("O27").Value = 75.68181818181819
$ws.Range("P27").Value = 44
$ws.Range("Q27").Value = 9
$ws.Range("R27").Value = 9
$ws.Range("S27").Value = 9

# Delete row 28 (only 27 data rows remain, dimension becomes A1:S27)
$ws.Rows.Item(28).Delete()
